# Edit script: rewrite the "El problema / La solución planteada / El alcance del
# proyecto" section of the document with the new expanded content, and drop
# the old "Especificaciones del prototipo" heading paragraph (superseded by a
# new closing paragraph about the module behaviour).

$d = $word.ActiveDocument

# The block we are replacing spans from the first empty paragraph right
# after the author list (paragraph 6) through the old "Especificaciones del
# prototipo" heading (paragraph 12) -- everything up to (but excluding) the
# three trailing empty paragraphs that close the document.
$pStart = $d.Paragraphs(6)
$pEnd = $d.Paragraphs(12)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Subttulo"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">El problema </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Subttulo"/>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/>
      <w:color w:val="auto"/>
      <w:spacing w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/>
      <w:color w:val="auto"/>
      <w:spacing w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>En Tumaco se localiza la mayor concentración de manglares del país, caracterizados por su riqueza natural y por sus valiosos recursos hidrobiológicos, siendo de gran relevancia para la región</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/>
      <w:color w:val="auto"/>
      <w:spacing w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> por lo que el ingreso a esta zona natural se encuentra </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/>
      <w:color w:val="auto"/>
      <w:spacing w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>restringido</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/>
      <w:color w:val="auto"/>
      <w:spacing w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> para evitar la contaminación del ecosistema o otras afectaciones que pueda tener la presencia del hombre en la zona. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Actualmente, </w:t>
  </w:r>
  <w:r>
    <w:t>la zona de manglares se encuentra</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> en alto riesgo de desecación para uso en construcción, para ampliar la zona urbana de Tumaco, por lo que es necesario el uso de sistemas de vigilancia y protección del área para evitar la incursión de personas con esta intención.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Subttulo"/>
  </w:pPr>
  <w:r>
    <w:t>La solución planteada</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Se propone</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> la implementación de un sistema de seguridad que sara aviso si se llega a detectar el ingreso de un intruso en una zona determinada </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">de la zona protegida, de preferencia a los accesos fluviales al manglar por medio </w:t>
  </w:r>
  <w:r>
    <w:t>del uso</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> de sensores ubicados </w:t>
  </w:r>
  <w:r>
    <w:t>en puntos estratégicos del manglar tales como desembocaduras de los afluentes a la zona, inicialmente.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Subttulo"/>
  </w:pPr>
  <w:r>
    <w:t>El alcance del proyecto</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">Se diseñará, modelará </w:t>
  </w:r>
  <w:r>
    <w:t>e</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> implementará un prototipo que consiste en </w:t>
  </w:r>
  <w:r>
    <w:t>un modulo emisor tipo TOF (</w:t>
  </w:r>
  <w:r>
    <w:t>Time-</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>of</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>-Flight</w:t>
  </w:r>
  <w:r>
    <w:t>), el cual detectará la distancia de un objeto fijo en el entorno por medio del tiempo que le toma al laser ir y volver al módulo.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t xml:space="preserve">El modulo al detectar un tiempo de retorno inconsistente con el registro base activara las alertas, una luminosa, otra sonora que se mantendrán activas hasta el momento que un operario las desactive </w:t>
  </w:r>
</w:p>
'@

$rng.InsertXML($newXml)

# InsertXML over a range that starts exactly on a paragraph boundary leaves
# that leading paragraph behind as an empty shell (its mark isn't consumed),
# so the freshly inserted paragraphs land just after it. Drop that leftover
# empty paragraph to land on the intended structure.
$d.Paragraphs(6).Range.Delete()

# Cosmetic: the "Default Paragraph Font" character style becomes referenced
# once the pasted-in runs above carry explicit run formatting, so Word
# un-hides it (drops <w:semiHidden/>) the next time the style sheet is
# written out.
try {
    $st = $d.Styles("Default Paragraph Font")
    $st.Hidden = $false
} catch {
}
